# "HC Services.xlsx" — add three new names to the "Main" tracking table:
#   - Universal Health (UHS), inserted right after LabCorp (row 20)
#   - Hapvida (HAPV3 BZ) and Chemed (CHE), appended at the bottom of the table

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Insert a new row above the current row 20 ("Catalent"); this pushes every
# row from Catalent down through Davita one row further down the sheet.
$ws.Rows.Item(20).Insert()

$ws.Range("A20").Value = "x"
$ws.Range("B20").Value = "Universal Health"
$ws.Range("C20").Value = "UHS"

# Append two brand-new rows at the bottom of the table.
$ws.Range("A43").Value = "x"
$ws.Range("B43").Value = "Hapvida"
$ws.Range("C43").Value = "HAPV3 BZ"

$ws.Range("A44").Value = "x"
$ws.Range("B44").Value = "Chemed"
$ws.Range("C44").Value = "CHE"

# The "Science 37" row (with the live formulas) shifted from row 33 to row
# 34 along with everything else, but its external hyperlink doesn't follow
# a plain row-insert automatically, so rebuild the sheet's hyperlinks at
# their (possibly new) locations.
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B4"), "CVS.xlsx") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "UNH.xlsx") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B34"), "SNCE.xlsx") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "ELV.xlsx") | Out-Null

# Re-adding a hyperlink resets the cell's style to a fresh (if equivalent)
# "Hyperlink" xf; re-apply the named style so the cells keep using the
# original style record instead of an accidental duplicate.
$ws.Range("B4").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("B34").Style = "Hyperlink"
$ws.Range("B5").Style = "Hyperlink"

# Restore the frozen-pane selection to reflect the newly-added last row.
$ws.Range("A43").Select() | Out-Null
